# Added serial run capture data.
#
# Inserts a new "Serial coWPAtty run time (milliseconds)" column (G) into
# the "Aggregate Tests" sheet (shifting the old "Dist Cow" column to H and
# "Result" to I), back-fills "N/A" for the pre-existing (distributed-run)
# rows, marks the two "no solution" test rows in red, and appends five new
# rows (7-11) capturing the serial coWPAtty run times.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert new column G; old G (Dist Cow) -> H, old H (Result) -> I ---
$ws.Columns("G:G").Insert()
$ws.Columns("G:G").ColumnWidth = 60.6

# Write the very first new row first so the shared-string table picks up
# "Serial" / "N/A" / "N/A - cmdline ..." / "Serial coWPAtty ..." in the same
# order as the original edit.
$ws.Range("A7").Value2 = "rbeede"
$ws.Range("B7").Value2 = "Serial"
$ws.Range("C7").Value2 = "linksys_FirstDictionary_!8zj39le"
$ws.Range("E7").Value2 = "N/A"
$ws.Range("D7").Value2 = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("F7").Value2 = "N/A"
$ws.Range("G7").Value2 = "N/A"
$ws.Range("H7").Value2 = 2
$ws.Range("I7").Value2 = "Correct - Password Found"

# --- New column header + "N/A" back-fill for the existing rows 2-6 ---
$ws.Range("G1").Value2 = "Serial coWPAtty run time (milliseconds)"
$ws.Range("G2").Value2 = "N/A"
$ws.Range("G3").Value2 = "N/A"
$ws.Range("G4").Value2 = "N/A"
$ws.Range("G5").Value2 = "N/A"
$ws.Range("G6").Value2 = "N/A"

# --- Highlight the "password not found" test-data rows in red text ---
$ws.Range("C5").Font.Color = 255
$ws.Range("C10").Font.Color = 255

# --- Row 8 ---
$ws.Range("A8").Value2 = "rbeede"
$ws.Range("B8").Value2 = "Serial"
$ws.Range("C8").Value2 = "linksys_MiddleDictionary_korrelie"
$ws.Range("D8").Value2 = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E8").Value2 = "N/A"
$ws.Range("F8").Value2 = "N/A"
$ws.Range("G8").Value2 = "N/A"
$ws.Range("H8").Value2 = 2981
$ws.Range("I8").Value2 = "Correct - Password Found"

# --- Row 9 ---
$ws.Range("A9").Value2 = "rbeede"
$ws.Range("B9").Value2 = "Serial"
$ws.Range("C9").Value2 = "linksys_LastDictionary_}ttringe"
$ws.Range("D9").Value2 = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E9").Value2 = "N/A"
$ws.Range("F9").Value2 = "N/A"
$ws.Range("G9").Value2 = "N/A"
$ws.Range("H9").Value2 = 5950
$ws.Range("I9").Value2 = "Correct - Password Found"

# --- Row 10 (no run time captured -> H10 left blank) ---
$ws.Range("A10").Value2 = "rbeede"
$ws.Range("B10").Value2 = "Serial"
$ws.Range("C10").Value2 = "linksys_NotInDictionary_UnknownPassword5763"
$ws.Range("D10").Value2 = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E10").Value2 = "N/A"
$ws.Range("F10").Value2 = "N/A"
$ws.Range("G10").Value2 = "N/A"
$ws.Range("I10").Value2 = "Correct - No Solution"

# --- Row 11 ---
$ws.Range("A11").Value2 = "rbeede"
$ws.Range("B11").Value2 = "Serial"
$ws.Range("C11").Value2 = "wireless_Test_invalid_capture"
$ws.Range("D11").Value2 = "N/A - cmdline = time cowpatty -d, -r, -s"
$ws.Range("E11").Value2 = "N/A"
$ws.Range("F11").Value2 = "N/A"
$ws.Range("G11").Value2 = "N/A"
$ws.Range("H11").Value2 = 3
$ws.Range("I11").Value2 = "Correct - No Solution"

# --- Update the view: scroll back to column A and select A11 ---
$ws.Activate()
$ws.Range("A1").Select()
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A11").Select()
